$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Mayo de 2020 a las 17:05"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1573777
$ws.Cells.Item(4, 3).Value = 3194
$ws.Cells.Item(4, 4).Value = 361419
$ws.Cells.Item(4, 5).Value = 1118653
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 172
$ws.Cells.Item(4, 8).Value = 93705

# Row 38
$ws.Cells.Item(38, 1).Value = "Rumania"
$ws.Cells.Item(38, 2).Value = 17387
$ws.Cells.Item(38, 3).Value = 196
$ws.Cells.Item(38, 4).Value = 10356
$ws.Cells.Item(38, 5).Value = 5887
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 7
$ws.Cells.Item(38, 8).Value = 1144

# Row 61
$ws.Cells.Item(61, 1).Value = "Moldavia"
$ws.Cells.Item(61, 2).Value = 6553
$ws.Cells.Item(61, 3).Value = 213
$ws.Cells.Item(61, 4).Value = 2953
$ws.Cells.Item(61, 5).Value = 3372
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 7
$ws.Cells.Item(61, 8).Value = 228

# Row 62
$ws.Cells.Item(62, 1).Value = "Finlandia"
$ws.Cells.Item(62, 2).Value = 6443
$ws.Cells.Item(62, 3).Value = 44
$ws.Cells.Item(62, 4).Value = 5000
$ws.Cells.Item(62, 5).Value = 1139
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 3
$ws.Cells.Item(62, 8).Value = 304

# Row 63
$ws.Cells.Item(63, 1).Value = "Nigeria"
$ws.Cells.Item(63, 2).Value = 6401
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 1734
$ws.Cells.Item(63, 5).Value = 4475
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 192

# Row 65
$ws.Cells.Item(65, 1).Value = "Oman"
$ws.Cells.Item(65, 2).Value = 6043
$ws.Cells.Item(65, 3).Value = 372
$ws.Cells.Item(65, 4).Value = 1661
$ws.Cells.Item(65, 5).Value = 4353
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 2
$ws.Cells.Item(65, 8).Value = 29

# Row 77
$ws.Cells.Item(77, 1).Value = "Grecia"
$ws.Cells.Item(77, 2).Value = 2850
$ws.Cells.Item(77, 3).Value = 10
$ws.Cells.Item(77, 4).Value = 1374
$ws.Cells.Item(77, 5).Value = 1310
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 166

# Row 123
$ws.Cells.Item(123, 1).Value = "San Marino"
$ws.Cells.Item(123, 2).Value = 656
$ws.Cells.Item(123, 3).Value = 1
$ws.Cells.Item(123, 4).Value = 220
$ws.Cells.Item(123, 5).Value = 395
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 41

# Row 135
$ws.Cells.Item(135, 1).Value = "Estado de Palestina"
$ws.Cells.Item(135, 2).Value = 398
$ws.Cells.Item(135, 3).Value = 7
$ws.Cells.Item(135, 4).Value = 346
$ws.Cells.Item(135, 5).Value = 50
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 2

# Row 139
$ws.Cells.Item(139, 1).Value = "Cabo Verde"
$ws.Cells.Item(139, 2).Value = 349
$ws.Cells.Item(139, 3).Value = 14
$ws.Cells.Item(139, 4).Value = 85
$ws.Cells.Item(139, 5).Value = 261
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 3

# Row 140
$ws.Cells.Item(140, 1).Value = "Togo"
$ws.Cells.Item(140, 2).Value = 338
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 107
$ws.Cells.Item(140, 5).Value = 219
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 12

# Row 150
$ws.Cells.Item(150, 1).Value = "Liberia"
$ws.Cells.Item(150, 2).Value = 238
$ws.Cells.Item(150, 3).Value = 5
$ws.Cells.Item(150, 4).Value = 128
$ws.Cells.Item(150, 5).Value = 87
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 23

# Row 157
$ws.Cells.Item(157, 1).Value = "Mozambique"
$ws.Cells.Item(157, 2).Value = 156
$ws.Cells.Item(157, 3).Value = 10
$ws.Cells.Item(157, 4).Value = 48
$ws.Cells.Item(157, 5).Value = 108
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 0

# Row 158
$ws.Cells.Item(158, 1).Value = "Guadalupe"
$ws.Cells.Item(158, 2).Value = 155
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 109
$ws.Cells.Item(158, 5).Value = 33
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 13

# Row 159
$ws.Cells.Item(159, 1).Value = "Gibraltar"
$ws.Cells.Item(159, 2).Value = 147
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 145
$ws.Cells.Item(159, 5).Value = 2
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

# Row 175
$ws.Cells.Item(175, 1).Value = "Malaui"
$ws.Cells.Item(175, 2).Value = 71
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(175, 4).Value = 27
$ws.Cells.Item(175, 5).Value = 41
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 3

# Row 183
$ws.Cells.Item(183, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(183, 2).Value = 40
$ws.Cells.Item(183, 3).Value = 1
$ws.Cells.Item(183, 4).Value = 33
$ws.Cells.Item(183, 5).Value = 4
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 3

# Row 184
$ws.Cells.Item(184, 1).Value = "Puerto Rico"
$ws.Cells.Item(184, 2).Value = 39
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 1
$ws.Cells.Item(184, 5).Value = 36
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 2

# Row 196
$ws.Cells.Item(196, 1).Value = "Santa Lucia"
$ws.Cells.Item(196, 2).Value = 18
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 18
$ws.Cells.Item(196, 5).Value = 0
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 0

# Row 197
$ws.Cells.Item(197, 1).Value = "Belice"
$ws.Cells.Item(197, 2).Value = 18
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 16
$ws.Cells.Item(197, 5).Value = 0
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 2

# Row 209
$ws.Cells.Item(209, 1).Value = "Montserrat"
$ws.Cells.Item(209, 2).Value = 11
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 10
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 1

# Row 210
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0
